$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "75.927.85"
$ws.Cells.Item(2, 5).Value = "  -0.21%  "

$ws.Cells.Item(3, 4).Value = "2.897.99"
$ws.Cells.Item(3, 5).Value = "  +1.51%  "

$ws.Cells.Item(4, 5).Value = "  +0.12%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "197.13"
$ws.Cells.Item(5, 5).Value = "  +1.66%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "589.51"
$ws.Cells.Item(6, 5).Value = "  -1.63%  "

$ws.Cells.Item(7, 5).Value = "  +0.06%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.544"
$ws.Cells.Item(8, 5).Value = "  -1.43%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.192"
$ws.Cells.Item(9, 5).Value = "  -1.18%  "

$ws.Cells.Item(10, 4).Value = "2.898.57"
$ws.Cells.Item(10, 5).Value = "  +1.70%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.443"
$ws.Cells.Item(11, 5).Value = "  +12.77%  "

$ws.Cells.Item(12, 5).Value = "  +0.29%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "4.85"
$ws.Cells.Item(13, 5).Value = "  -1.10%  "

$ws.Cells.Item(14, 4).Value = "3.435.86"
$ws.Cells.Item(14, 5).Value = "  +2.18%  "

$ws.Cells.Item(15, 4).Value = "75.886.05"
$ws.Cells.Item(15, 5).Value = "  -0.01%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "27.67"
$ws.Cells.Item(16, 5).Value = "  +0.48%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.0000185"
$ws.Cells.Item(17, 5).Value = "  -2.24%  "

$ws.Cells.Item(18, 4).Value = "2.896.36"
$ws.Cells.Item(18, 5).Value = "  +1.78%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "13.04"
$ws.Cells.Item(19, 5).Value = "  +4.18%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "8.60"
$ws.Cells.Item(20, 5).Value = "  -6.10%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "365.76"
$ws.Cells.Item(21, 5).Value = "  -4.45%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.27"
$ws.Cells.Item(22, 5).Value = "  +3.00%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "2.22"
$ws.Cells.Item(23, 5).Value = "  -4.75%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "71.61"
$ws.Cells.Item(24, 5).Value = "  -0.60%  "

$ws.Cells.Item(25, 5).Value = "  -0.18%  "

$ws.Cells.Item(26, 4).Value = "3.060.16"
$ws.Cells.Item(26, 5).Value = "  +2.32%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "4.19"
$ws.Cells.Item(27, 5).Value = "  -0.80%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "9.49"
$ws.Cells.Item(28, 5).Value = "  -2.64%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.0000104"
$ws.Cells.Item(29, 5).Value = "  -0.90%  "

$ws.Cells.Item(30, 5).Value = "  +0.31%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "7.95"
$ws.Cells.Item(31, 5).Value = "  +3.24%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.35"
$ws.Cells.Item(32, 5).Value = "  -5.05%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "490.38"
$ws.Cells.Item(33, 5).Value = "  -5.56%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.80"
$ws.Cells.Item(34, 5).Value = "  -0.98%  "

$ws.Cells.Item(35, 5).Value = "  +0.11%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "165.11"
$ws.Cells.Item(36, 5).Value = "  -0.79%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "19.90"
$ws.Cells.Item(37, 5).Value = "  -0.68%  "

$ws.Cells.Item(38, 2).Value = "Cronos"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.109"
$ws.Cells.Item(38, 5).Value = "  +22.41%  "

$ws.Cells.Item(39, 2).Value = "PolygonEcosystemToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.380"
$ws.Cells.Item(39, 5).Value = "  +10.31%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "19.67"
$ws.Cells.Item(40, 5).Value = "  +1.26%  "

$ws.Cells.Item(41, 5).Value = "  +0.00%  "

$ws.Cells.Item(42, 5).Value = "  -9.10%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "178.32"
$ws.Cells.Item(43, 5).Value = "  -4.44%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "4.81"
$ws.Cells.Item(44, 5).Value = "  -5.72%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.62"
$ws.Cells.Item(45, 5).Value = "  -3.78%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "40.05"
$ws.Cells.Item(46, 5).Value = "  -0.61%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.16"
$ws.Cells.Item(47, 5).Value = "  -6.20%  "

$ws.Cells.Item(48, 2).Value = "Filecoin"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "3.81"
$ws.Cells.Item(48, 5).Value = "  +1.20%  "

$ws.Cells.Item(49, 2).Value = "ARBITRUM"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.571"
$ws.Cells.Item(49, 5).Value = "  -0.71%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.21"
$ws.Cells.Item(50, 5).Value = "  -6.84%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "22.07"
$ws.Cells.Item(51, 5).Value = "  +2.67%  "
